# Update team-specific time-matrix probabilities on Sheet1 (Buffalo_B).
# Only the numeric transition-probability cells that changed between the
# before/after commit are touched; zero-valued and label cells are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2244897959183673
$ws.Range("C2").Value = 0.5612244897959183
$ws.Range("J2").Value = 0.01700680272108844
$ws.Range("P2").Value = 0.119047619047619
$ws.Range("S2").Value = 0.07823129251700681
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.02352941176470588
$ws.Range("J3").Value = 0.02352941176470588
$ws.Range("P3").Value = 0.7529411764705882
$ws.Range("S3").Value = 0.1882352941176471
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.813953488372093
$ws.Range("S4").Value = 0.1627906976744186
$ws.Range("B6").Value = 0.0639269406392694
$ws.Range("D6").Value = 0.0091324200913242
$ws.Range("F6").Value = 0.0593607305936073
$ws.Range("J6").Value = 0.2602739726027397
$ws.Range("O6").Value = 0.0091324200913242
$ws.Range("Q6").Value = 0.228310502283105
$ws.Range("R6").Value = 0.0730593607305936
$ws.Range("S6").Value = 0.2968036529680365
$ws.Range("B7").Value = 0.09467455621301775
$ws.Range("D7").Value = 0.02366863905325444
$ws.Range("F7").Value = 0.03550295857988166
$ws.Range("J7").Value = 0.136094674556213
$ws.Range("O7").Value = 0.01775147928994083
$ws.Range("Q7").Value = 0.242603550295858
$ws.Range("R7").Value = 0.08875739644970414
$ws.Range("S7").Value = 0.3609467455621302
$ws.Range("B8").Value = 0.08951965065502183
$ws.Range("D8").Value = 0.01310043668122271
$ws.Range("F8").Value = 0.06331877729257641
$ws.Range("J8").Value = 0.1004366812227074
$ws.Range("O8").Value = 0.01746724890829694
$ws.Range("Q8").Value = 0.2510917030567685
$ws.Range("R8").Value = 0.07641921397379912
$ws.Range("S8").Value = 0.388646288209607
$ws.Range("B9").Value = 0.1067415730337079
$ws.Range("D9").Value = 0.02247191011235955
$ws.Range("F9").Value = 0.03370786516853932
$ws.Range("J9").Value = 0.07865168539325842
$ws.Range("O9").Value = 0.02808988764044944
$ws.Range("Q9").Value = 0.2134831460674157
$ws.Range("R9").Value = 0.1123595505617977
$ws.Range("S9").Value = 0.4044943820224719
$ws.Range("B10").Value = 0.0974910394265233
$ws.Range("D10").Value = 0.01935483870967742
$ws.Range("E10").Value = 0.0007168458781362007
$ws.Range("F10").Value = 0.06594982078853047
$ws.Range("J10").Value = 0.1290322580645161
$ws.Range("O10").Value = 0.01075268817204301
$ws.Range("Q10").Value = 0.2365591397849462
$ws.Range("R10").Value = 0.08960573476702509
$ws.Range("S10").Value = 0.3505376344086021
$ws.Range("G11").Value = 0.1147540983606557
$ws.Range("J11").Value = 0.0860655737704918
$ws.Range("K11").Value = 0.1762295081967213
$ws.Range("L11").Value = 0.6065573770491803
$ws.Range("S11").Value = 0.01639344262295082
$ws.Range("G12").Value = 0.7197452229299363
$ws.Range("J12").Value = 0.1847133757961783
$ws.Range("K12").Value = 0.006369426751592357
$ws.Range("L12").Value = 0.05732484076433121
$ws.Range("S12").Value = 0.03184713375796178
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.5918367346938775
$ws.Range("J13").Value = 0.3061224489795918
$ws.Range("S13").Value = 0.08163265306122448
$ws.Range("F15").Value = 0.02985074626865672
$ws.Range("H15").Value = 0.1641791044776119
$ws.Range("I15").Value = 0.03482587064676617
$ws.Range("J15").Value = 0.417910447761194
$ws.Range("K15").Value = 0.05970149253731343
$ws.Range("M15").Value = 0.009950248756218905
$ws.Range("O15").Value = 0.05970149253731343
$ws.Range("S15").Value = 0.2238805970149254
$ws.Range("F16").Value = 0.0053475935828877
$ws.Range("H16").Value = 0.160427807486631
$ws.Range("I16").Value = 0.06417112299465241
$ws.Range("J16").Value = 0.4973262032085561
$ws.Range("K16").Value = 0.08021390374331551
$ws.Range("M16").Value = 0.0160427807486631
$ws.Range("O16").Value = 0.0481283422459893
$ws.Range("S16").Value = 0.1283422459893048
$ws.Range("F17").Value = 0.02086956521739131
$ws.Range("H17").Value = 0.1443478260869565
$ws.Range("I17").Value = 0.08695652173913043
$ws.Range("J17").Value = 0.4695652173913044
$ws.Range("K17").Value = 0.07130434782608695
$ws.Range("M17").Value = 0.01913043478260869
$ws.Range("O17").Value = 0.05565217391304348
$ws.Range("S17").Value = 0.1321739130434783
$ws.Range("F18").Value = 0.004739336492890996
$ws.Range("H18").Value = 0.1753554502369668
$ws.Range("I18").Value = 0.09004739336492891
$ws.Range("J18").Value = 0.4454976303317535
$ws.Range("K18").Value = 0.1184834123222749
$ws.Range("M18").Value = 0.01421800947867299
$ws.Range("O18").Value = 0.04265402843601896
$ws.Range("S18").Value = 0.1090047393364929
$ws.Range("F19").Value = 0.01879084967320261
$ws.Range("H19").Value = 0.2295751633986928
$ws.Range("I19").Value = 0.0727124183006536
$ws.Range("J19").Value = 0.3897058823529412
$ws.Range("K19").Value = 0.0857843137254902
$ws.Range("M19").Value = 0.02532679738562092
$ws.Range("O19").Value = 0.06535947712418301
$ws.Range("S19").Value = 0.1127450980392157
